# Update crypto price/volume table cells per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.723.85'
$ws.Range('E2').Value = '  +2.51%  '
$ws.Range('D3').Value = '2.435.19'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.37'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.01%  '
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.93'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.30'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.01'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('D15').Value = '2.814.72'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').Value = '2.442.41'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.839'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '45.602.52'
$ws.Range('E18').Value = '  +2.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.40'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').Value = '0.0₃0929'
$ws.Range('E21').Value = '  +2.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.88%  '
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '246.49'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.42%  '
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.18%  '
$ws.Range('E28').Value = '  -4.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.67'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.53'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '49.37'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.129'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.15'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0760'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.53'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.91'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '127.80'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('E42').Value = '  +1.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.83'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0292'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').Value = '1.957.78'
$ws.Range('E45').Value = '  +0.46%  '
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('E48').Value = '  +7.84%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.14'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '77.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.84'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.04%  '
